$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-23 Monday" "2024-12-24 Tuesday"

Replace-Text "433÷9=48, 1" "455÷8=56, 7"
Replace-Text "674÷4=168, 2" "610÷9=67, 7"
Replace-Text "366÷9=40, 6" "324÷7=46, 2"
Replace-Text "844÷5=168, 4" "997÷4=249, 1"
Replace-Text "866÷5=173, 1" "421÷5=84, 1"

Replace-Text "580÷7=82, 6" "981÷2=490, 1"
Replace-Text "746÷2=373, 0" "213÷3=71, 0"
Replace-Text "402÷3=134, 0" "512÷5=102, 2"
Replace-Text "405÷3=135, 0" "549÷3=183, 0"
Replace-Text "768÷9=85, 3" "325÷2=162, 1"

Replace-Text "119÷2=59, 1" "441÷9=49, 0"
Replace-Text "476÷5=95, 1" "175÷5=35, 0"
Replace-Text "752÷7=107, 3" "823÷2=411, 1"
Replace-Text "905÷9=100, 5" "611÷2=305, 1"
Replace-Text "904÷9=100, 4" "755÷9=83, 8"

Replace-Text "384÷4=96, 0" "321÷2=160, 1"
Replace-Text "259÷6=43, 1" "645÷9=71, 6"
Replace-Text "636÷5=127, 1" "504÷6=84, 0"
Replace-Text "889÷6=148, 1" "220÷7=31, 3"
Replace-Text "945÷6=157, 3" "894÷5=178, 4"

Replace-Text "519÷4=129, 3" "387÷8=48, 3"
Replace-Text "321÷6=53, 3" "311÷3=103, 2"
Replace-Text "810÷3=270, 0" "471÷9=52, 3"
Replace-Text "529÷9=58, 7" "941÷9=104, 5"
Replace-Text "306÷2=153, 0" "637÷8=79, 5"

Write-Output "Replacements applied"
